$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 1810.6666
$ws.Cells.Item(111, 9).Value = 1500
$ws.Cells.Item(111, 11).Value = 4500
$ws.Cells.Item(111, 13).Value = -1433

$ws.Cells.Item(113, 8).Value = 3589.3667
$ws.Cells.Item(113, 9).Value = 3060.4167
$ws.Cells.Item(113, 10).Value = 3942
$ws.Cells.Item(113, 11).Value = 3060.4167
$ws.Cells.Item(113, 12).Value = 3942
$ws.Cells.Item(113, 13).Value = 193.5832999999998
$ws.Cells.Item(113, 14).Value = -10450

$ws.Cells.Item(127, 8).Value = 332.4375
$ws.Cells.Item(127, 9).Value = 332.4375
$ws.Cells.Item(127, 11).Value = 997.3125
$ws.Cells.Item(127, 13).Value = 3962.6875

$ws.Cells.Item(137, 8).Value = 340190.66
$ws.Cells.Item(137, 9).Value = 2838.5356
$ws.Cells.Item(137, 11).Value = 8515.606800000001
$ws.Cells.Item(137, 13).Value = -5965.606800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 11906735
$ws.Cells.Item(74, 9).Value = 7937694.5
$ws.Cells.Item(74, 10).Value = 23813856
$ws.Cells.Item(74, 11).Value = 7937694.5
$ws.Cells.Item(74, 12).Value = 23813856
$ws.Cells.Item(74, 13).Value = -7936820.5
$ws.Cells.Item(74, 14).Value = -23815604

$ws.Cells.Item(77, 8).Value = 11906735
$ws.Cells.Item(77, 9).Value = 7937694.5
$ws.Cells.Item(77, 10).Value = 23813856
$ws.Cells.Item(77, 11).Value = 39688472.5
$ws.Cells.Item(77, 12).Value = 119069280
$ws.Cells.Item(77, 13).Value = -39684104.5
$ws.Cells.Item(77, 14).Value = -119078016

$ws.Cells.Item(110, 8).Value = 2963.4443
$ws.Cells.Item(110, 9).Value = 2808.875
$ws.Cells.Item(110, 10).Value = 4200
$ws.Cells.Item(110, 11).Value = 2808.875
$ws.Cells.Item(110, 12).Value = 4200
$ws.Cells.Item(110, 13).Value = -763.875
$ws.Cells.Item(110, 14).Value = -8290

$ws.Cells.Item(122, 8).Value = 3709.6667
$ws.Cells.Item(122, 9).Value = 5900
$ws.Cells.Item(122, 10).Value = 2614.5
$ws.Cells.Item(122, 11).Value = 17700
$ws.Cells.Item(122, 12).Value = 7843.5
$ws.Cells.Item(122, 13).Value = -15250
$ws.Cells.Item(122, 14).Value = -12743.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 4860
$ws.Cells.Item(99, 9).Value = 4860
$ws.Cells.Item(99, 11).Value = 4860
$ws.Cells.Item(99, 13).Value = -3362

$ws.Cells.Item(134, 8).Value = 19955.678
$ws.Cells.Item(134, 9).Value = 25108.953
$ws.Cells.Item(134, 10).Value = 6106.25
$ws.Cells.Item(134, 11).Value = 75326.859
$ws.Cells.Item(134, 12).Value = 18318.75
$ws.Cells.Item(134, 13).Value = -72791.859
$ws.Cells.Item(134, 14).Value = -23388.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3016.614
$ws.Cells.Item(31, 9).Value = 2356.1765
$ws.Cells.Item(31, 10).Value = 3992.913
$ws.Cells.Item(31, 11).Value = 2356.1765
$ws.Cells.Item(31, 12).Value = 3992.913
$ws.Cells.Item(31, 13).Value = -2061.1765
$ws.Cells.Item(31, 14).Value = -4582.913

$ws.Cells.Item(34, 8).Value = 3016.614
$ws.Cells.Item(34, 9).Value = 2356.1765
$ws.Cells.Item(34, 10).Value = 3992.913
$ws.Cells.Item(34, 11).Value = 2356.1765
$ws.Cells.Item(34, 12).Value = 3992.913
$ws.Cells.Item(34, 13).Value = -2154.1765
$ws.Cells.Item(34, 14).Value = -4396.913

$ws.Cells.Item(59, 8).Value = 63333.332
$ws.Cells.Item(59, 10).Value = 75000
$ws.Cells.Item(59, 12).Value = 75000
$ws.Cells.Item(59, 14).Value = -77290

$ws.Cells.Item(60, 8).Value = 29333.834
$ws.Cells.Item(60, 10).Value = 29333.834
$ws.Cells.Item(60, 12).Value = 29333.834
$ws.Cells.Item(60, 14).Value = -30355.834

$ws.Cells.Item(141, 8).Value = 25558.77
$ws.Cells.Item(141, 10).Value = 24355.334
$ws.Cells.Item(141, 12).Value = 24355.334
$ws.Cells.Item(141, 14).Value = -34715.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 669.25
$ws.Cells.Item(113, 9).Value = 836
$ws.Cells.Item(113, 10).Value = 569.2
$ws.Cells.Item(113, 11).Value = 2508
$ws.Cells.Item(113, 12).Value = 1707.6
$ws.Cells.Item(113, 13).Value = -338
$ws.Cells.Item(113, 14).Value = -6047.6

$ws.Cells.Item(124, 8).Value = 1839.375
$ws.Cells.Item(124, 10).Value = 2233.3333
$ws.Cells.Item(124, 12).Value = 6699.999899999999
$ws.Cells.Item(124, 14).Value = -16519.9999

$ws.Cells.Item(131, 8).Value = 1407.5469
$ws.Cells.Item(131, 10).Value = 1577.7646
$ws.Cells.Item(131, 12).Value = 4733.293799999999
$ws.Cells.Item(131, 14).Value = -14813.2938

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 9421.857
$ws.Cells.Item(57, 10).Value = 10599.6
$ws.Cells.Item(57, 12).Value = 10599.6
$ws.Cells.Item(57, 14).Value = -12239.6

$ws.Cells.Item(70, 8).Value = 14333.25
$ws.Cells.Item(70, 9).Value = 3966.6667
$ws.Cells.Item(70, 10).Value = 24699.834
$ws.Cells.Item(70, 11).Value = 3966.6667
$ws.Cells.Item(70, 12).Value = 24699.834
$ws.Cells.Item(70, 13).Value = -3696.6667
$ws.Cells.Item(70, 14).Value = -25239.834

$ws.Cells.Item(73, 8).Value = 14333.25
$ws.Cells.Item(73, 9).Value = 3966.6667
$ws.Cells.Item(73, 10).Value = 24699.834
$ws.Cells.Item(73, 11).Value = 3966.6667
$ws.Cells.Item(73, 12).Value = 24699.834
$ws.Cells.Item(73, 13).Value = -3030.6667
$ws.Cells.Item(73, 14).Value = -26571.834

$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).ClearContents()

$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).ClearContents()

$ws.Cells.Item(102, 8).Value = 3834.8333
$ws.Cells.Item(102, 9).Value = 4159.7896
$ws.Cells.Item(102, 10).Value = 2600
$ws.Cells.Item(102, 11).Value = 4159.7896
$ws.Cells.Item(102, 12).Value = 2600
$ws.Cells.Item(102, 13).Value = -2537.7896
$ws.Cells.Item(102, 14).Value = -5844

$ws.Cells.Item(122, 8).Value = 1886.8572
$ws.Cells.Item(122, 9).Value = 1714.2858
$ws.Cells.Item(122, 10).Value = 2059.4285
$ws.Cells.Item(122, 11).Value = 5142.857400000001
$ws.Cells.Item(122, 12).Value = 6178.2855
$ws.Cells.Item(122, 13).Value = -2692.857400000001
$ws.Cells.Item(122, 14).Value = -11078.2855

$ws.Cells.Item(126, 8).Value = 2887.8948
$ws.Cells.Item(126, 9).Value = 2928.2727
$ws.Cells.Item(126, 10).Value = 2832.375
$ws.Cells.Item(126, 11).Value = 8784.8181
$ws.Cells.Item(126, 12).Value = 8497.125
$ws.Cells.Item(126, 13).Value = -6314.8181
$ws.Cells.Item(126, 14).Value = -13437.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2617.5
$ws.Cells.Item(7, 9).Value = 2541.6667
$ws.Cells.Item(7, 10).Value = 2754
$ws.Cells.Item(7, 11).Value = 2541.6667
$ws.Cells.Item(7, 12).Value = 2754
$ws.Cells.Item(7, 13).Value = -2429.6667
$ws.Cells.Item(7, 14).Value = -2978

$ws.Cells.Item(40, 8).Value = 2567.4666
$ws.Cells.Item(40, 9).Value = 2376
$ws.Cells.Item(40, 10).Value = 3333.3333
$ws.Cells.Item(40, 11).Value = 2376
$ws.Cells.Item(40, 12).Value = 3333.3333
$ws.Cells.Item(40, 13).Value = -2240
$ws.Cells.Item(40, 14).Value = -3605.3333

$ws.Cells.Item(61, 8).Value = 2139.5715
$ws.Cells.Item(61, 9).Value = 2139.5715
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 2139.5715
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -1937.5715
$ws.Cells.Item(61, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 2139.5715
$ws.Cells.Item(113, 9).Value = 2139.5715
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 2139.5715
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 30.42849999999999
$ws.Cells.Item(113, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 4110.6665
$ws.Cells.Item(122, 9).Value = 3791.125
$ws.Cells.Item(122, 10).Value = 4749.75
$ws.Cells.Item(122, 11).Value = 11373.375
$ws.Cells.Item(122, 12).Value = 14249.25
$ws.Cells.Item(122, 13).Value = -8923.375
$ws.Cells.Item(122, 14).Value = -19149.25

$ws.Cells.Item(126, 8).Value = 2617.5
$ws.Cells.Item(126, 9).Value = 2541.6667
$ws.Cells.Item(126, 10).Value = 2754
$ws.Cells.Item(126, 11).Value = 7625.000100000001
$ws.Cells.Item(126, 12).Value = 8262
$ws.Cells.Item(126, 13).Value = -5155.000100000001
$ws.Cells.Item(126, 14).Value = -13202

$ws.Cells.Item(136, 8).Value = 4109.0215
$ws.Cells.Item(136, 9).Value = 2080.4285
$ws.Cells.Item(136, 10).Value = 10025.75
$ws.Cells.Item(136, 11).Value = 6241.2855
$ws.Cells.Item(136, 12).Value = 30077.25
$ws.Cells.Item(136, 13).Value = -3691.2855
$ws.Cells.Item(136, 14).Value = -35177.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 57739.277
$ws.Cells.Item(122, 9).Value = 79239
$ws.Cells.Item(122, 10).Value = 1840
$ws.Cells.Item(122, 11).Value = 237717
$ws.Cells.Item(122, 12).Value = 5520
$ws.Cells.Item(122, 13).Value = -235267
$ws.Cells.Item(122, 14).Value = -10420

$ws.Cells.Item(136, 8).Value = 3082.4783
$ws.Cells.Item(136, 9).Value = 1926.1333
$ws.Cells.Item(136, 10).Value = 5250.625
$ws.Cells.Item(136, 11).Value = 5778.3999
$ws.Cells.Item(136, 12).Value = 15751.875
$ws.Cells.Item(136, 13).Value = -3228.3999
$ws.Cells.Item(136, 14).Value = -20851.875
